# Add 2022-Q4 data:
#  - insert a new "2022-Q4" worksheet (fund-holder detail) right before "2021-Q4"
#  - insert a corresponding summary row at the top of the "总计" sheet's data

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q4" worksheet before the existing "2021-Q4" sheet.
# ---------------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item("2021-Q4")
$new = $wb.Worksheets.Add($refSheet)
$new.Name = "2022-Q4"

# Copy header-row / index-column formatting from the "2021-Q4" sheet so the
# new sheet matches the look of its siblings (bold header, bordered cells).
$refSheet.Range("A1:H1").Copy()
$new.Range("A1:H1").PasteSpecial(-4122)
$refSheet.Range("A2:A3").Copy()
$new.Range("A2:A10").PasteSpecial(-4122)

# Header row
$new.Cells.Item(1,2).Value = "基金代码"
$new.Cells.Item(1,3).Value = "基金名称"
$new.Cells.Item(1,4).Value = "基金规模"
$new.Cells.Item(1,5).Value = "股票总仓位"
$new.Cells.Item(1,6).Value = "仓位占比"
$new.Cells.Item(1,7).Value = "持有市值(亿元)"
$new.Cells.Item(1,8).Value = "仓位排名"

# Fund rows: Index(A,n) Code(B,text) Name(C,text) Scale(D,text) Position(E,text) Ratio(F,text) MktValue(G,text) Rank(H,n)
$fundRows = @(
    @(0, "002291", "诺安安鑫灵活配置混合",               "2.40", "91.19", "4.22", "0.1013", 4),
    @(1, "002067", "诺安精选回报灵活配置混合",           "1.27", "50.95", "3.78", "0.0480", 3),
    @(2, "006167", "德邦乐享生活混合A",                   "1.01", "92.80", "2.91", "0.0294", 10),
    @(3, "006168", "德邦乐享生活混合C",                   "0.57", "92.80", "2.91", "0.0166", 10),
    @(4, "562520", "华夏中证智选1000成长创新策略ETF",     "0.38", "96.24", "1.00", "0.0038", 4),
    @(5, "012415", "德邦上证 G60 创新综合指数增强A",      "0.09", "91.89", "3.87", "0.0035", 4),
    @(6, "001412", "德邦鑫星价值灵活配置混合A",           "0.13", "35.79", "2.22", "0.0029", 5),
    @(7, "002112", "德邦鑫星价值灵活配置混合C",           "0.02", "35.79", "2.22", "0.0004", 5),
    @(8, "012416", "德邦上证 G60 创新综合指数增强C",      "0.01", "91.89", "3.87", "0.0004", 4)
)

for ($i = 0; $i -lt $fundRows.Count; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]

    $new.Cells.Item($r,1).Value = $row[0]

    $new.Cells.Item($r,2).NumberFormat = "@"
    $new.Cells.Item($r,2).Value = $row[1]

    $new.Cells.Item($r,3).NumberFormat = "@"
    $new.Cells.Item($r,3).Value = $row[2]

    $new.Cells.Item($r,4).NumberFormat = "@"
    $new.Cells.Item($r,4).Value = $row[3]

    $new.Cells.Item($r,5).NumberFormat = "@"
    $new.Cells.Item($r,5).Value = $row[4]

    $new.Cells.Item($r,6).NumberFormat = "@"
    $new.Cells.Item($r,6).Value = $row[5]

    $new.Cells.Item($r,7).NumberFormat = "@"
    $new.Cells.Item($r,7).Value = $row[6]

    $new.Cells.Item($r,8).Value = $row[7]
}

$new.Range("A1").Select()

# ---------------------------------------------------------------------------
# 2) Insert the "2022-Q4" summary row at the top of the "总计" sheet's data
#    (row 2), pushing the existing quarters down and renumbering the index
#    column (A) so it keeps counting 0,1,2,3,4,5.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(2).Insert()
# Row 3 (the old row 2, "2021-Q4") kept its original formatting; clone it
# into the freshly inserted row 2.
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)

$totalRows = @(
    @("2022-Q4", 9, 0.21),
    @("2021-Q4", 2, 0.23),
    @("2021-Q3", 1, 0.13),
    @("2021-Q2", 1, 0.06),
    @("2021-Q1", 2, 0),
    @("2020-Q4", 9, 1.47)
)

for ($i = 0; $i -lt $totalRows.Count; $i++) {
    $r = $i + 2
    $row = $totalRows[$i]
    $total.Cells.Item($r,1).Value = $i
    $total.Cells.Item($r,2).Value = $row[0]
    $total.Cells.Item($r,3).Value = $row[1]
    $total.Cells.Item($r,4).Value = $row[2]
}

$total.Range("A1").Select()
